$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in existing shared string used at A9 ("reques" -> "request")
$ws.Range("A9").Value = "Avance de los request de la web cliente."

# Fill in row 10 with the new log entry ("FInalización de los request de la web cliente.")
$ws.Range("A10").Value = "FInalización de los request de la web cliente."

$ws.Range("B10").NumberFormat = "m/d/yy h:mm"
$ws.Range("B10").Value = 42879.833333333336

$ws.Range("C10").NumberFormat = "m/d/yy h:mm"
$ws.Range("C10").Value = 42880.166666666664

$ws.Range("D10").NumberFormat = "h:mm"
$ws.Range("D10").Value = 0.33333333333333331

# Update the view: scroll so row 4 is at top, and move the selection to A10
$ws.Range("A10").Select()
$excel.ActiveWindow.ScrollRow = 4
